# Update workbook to the new data snapshot.
# The update re-sequences several fixture rows (their B..AB content, i.e.
# everything except the running index in column A, is swapped between two
# adjacent rows that share the same kick-off date) and removes the final
# two fixtures of the sheet (rows for ids 143 and 144).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of worksheet rows whose B:AB content (match id, teams, odds, ...)
# needs to be swapped with one another.
$pairs = @(
    @(16, 17),
    @(69, 70),
    @(86, 87),
    @(117, 118),
    @(125, 126),
    @(143, 144)
)

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    $range1 = $ws.Range("B$row1`:AB$row1")
    $range2 = $ws.Range("B$row2`:AB$row2")

    $val1 = $range1.Value()
    $val2 = $range2.Value()

    $range1.Value = $val2
    $range2.Value = $val1
}

# Remove the last two fixture rows (previously rows 145 and 146) entirely,
# shifting the remaining rows up. Delete the higher-numbered row first.
$ws.Rows.Item(146).Delete()
$ws.Rows.Item(145).Delete()
